# AcompReq.xlsx update
# Rows 84-87: fill in OF_CDG (F), OF_DATA (G), ITEM_PRCUNTPED (L), PRCTTL_INSUMO (M),
#             FORNECEDOR_CDG (N) and FORNECEDOR_DESC (O) -- previously blank/zero.
# Row 169/170: new "E.02.0041" order row fills in with real data and takes the
#             material code/description that used to sit on row 170
#             ("E.01.0114" / COMPACTADOR...); row 170 in turn is re-labelled
#             back to the material that row 169 used to carry
#             ("C.04.0100" / DESINFETANTE 5 L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumberCell($addr, $formatSourceAddr, $value) {
    # Copy number format from an existing formatted cell so we reuse the
    # same style index instead of minting a new one, then write the value.
    $ws.Range($formatSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $formatSourceAddr, $text) {
    # Force text storage (keep leading zeros / leading spaces) without
    # picking up a "@" custom-format style: mark the cell as Text, assign
    # the value, then restore the original (General) format from a
    # same-styled reference cell so the saved style index matches.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($formatSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------- Row 84 ----------
Set-NumberCell "F84" "F2" 80115
Set-NumberCell "G84" "G2" 45952
$ws.Range("L84").Value = 317.68
$ws.Range("M84").Value = 317.68
Set-TextCell "N84" "N2" "00000000008508"
Set-TextCell "O84" "O2" " MASTERPOL"

# ---------- Row 85 ----------
Set-NumberCell "F85" "F2" 80115
Set-NumberCell "G85" "G2" 45952
$ws.Range("L85").Value = 500.34592
$ws.Range("M85").Value = 27018.67968
Set-TextCell "N85" "N2" "00000000008508"
Set-TextCell "O85" "O2" " MASTERPOL"

# ---------- Row 86 ----------
Set-NumberCell "F86" "F2" 80115
Set-NumberCell "G86" "G2" 45952
$ws.Range("L86").Value = 672.46
$ws.Range("M86").Value = 672.46
Set-TextCell "N86" "N2" "00000000008508"
Set-TextCell "O86" "O2" " MASTERPOL"

# ---------- Row 87 ----------
Set-NumberCell "F87" "F2" 80115
Set-NumberCell "G87" "G2" 45952
$ws.Range("L87").Value = 655.46
$ws.Range("M87").Value = 26218.4
Set-TextCell "N87" "N2" "00000000008508"
Set-TextCell "O87" "O2" " MASTERPOL"

# ---------- Row 169 ----------
Set-NumberCell "F169" "F2" 80116
Set-NumberCell "G169" "G2" 45952
$ws.Range("H169").Value = "E.01.0114"
$ws.Range("I169").Value = "COMPACTADOR DE SOLO  À PERCUSSÃO, GASOLINA"
$ws.Range("L169").Value = 325
$ws.Range("M169").Value = 325
Set-TextCell "N169" "N2" "00000000007636"
Set-TextCell "O169" "O2" "ABRIL LOC"

# ---------- Row 170 ----------
$ws.Range("H170").Value = "C.04.0100"
$ws.Range("I170").Value = "DESINFETANTE 5 L"
